{"js": "const body = context.document.body;\n\n// --- Change 1: the momento badge is currently split across three runs\n//     (\"{\", \"momento\", \"}\") that all share identical formatting. Replacing\n//     the whole paragraph's text collapses them into a single run whose\n//     text reads \"{momento}\" (formatting of the first run is preserved). ---\nconst paras = body.paragraphs;\nparas.load(\"items/style\");\nawait context.sync();\n\nlet momentoPara = null;\nfor (const p of paras.items) {\n  if (p.style === \"mgmomento\") {\n    momentoPara = p;\n    break;\n  }\n}\nif (momentoPara) {\n  momentoPara.insertText(\"{momento}\", \"Replace\");\n}\n\n// --- Change 2: turn the lone \"{cifra}\" placeholder paragraph into a\n//     \"{#cifra}{.}\" ... \"{/cifra}\" loop block (two paragraphs). ---\nconst hits = body.search(\"{cifra}\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  const target = hits.items[0];\n  const cifraPara = target.paragraphs.getFirst();\n  target.insertText(\"{#cifra}{.}\", \"Replace\");\n  cifraPara.insertParagraph(\"{/cifra}\", \"After\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: the momento badge is currently split across three runs\n#     (\"{\", \"momento\", \"}\") that all share identical formatting. Collapse\n#     them into a single run whose text reads \"{momento}\". ---\n$pMomento = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"mgmomento\") {\n        $pMomento = $p\n        break\n    }\n}\nif ($pMomento -ne $null) {\n    # Range over the paragraph text only (exclude the trailing paragraph mark).\n    $rngMomento = $d.Range($pMomento.Range.Start, $pMomento.Range.End - 1)\n    # Word treats a same-value assignment as a no-op and leaves the run split\n    # untouched, so nudge the text through a throwaway value first to force a\n    # real rewrite (which merges the runs using the first run's formatting),\n    # then set the final value.\n    $rngMomento.Text = $rngMomento.Text + [char]1\n    $rngMomento = $d.Range($pMomento.Range.Start, $pMomento.Range.End - 1)\n    $rngMomento.Text = \"{momento}\"\n}\n\n# --- Change 2: turn the lone \"{cifra}\" placeholder paragraph into a\n#     \"{#cifra}{.}\" ... \"{/cifra}\" loop block (two paragraphs). ---\n$rngCifra = $d.Content\n$found = $rngCifra.Find.Execute(\"{cifra}\")\nif ($found) {\n    $cifraParaIndex = $rngCifra.Paragraphs.First.Index\n    $rngCifra.Text = \"{#cifra}{.}\"\n    $rngCifra.Collapse(0)\n    $rngCifra.InsertParagraphAfter()\n\n    $newPara = $d.Paragraphs.Item($cifraParaIndex + 1)\n    $newPara.Range.Text = \"{/cifra}\"\n}\n"}
